$d = $word.ActiveDocument

# Locate the "ÇOK-ÇOK" paragraph; the two new paragraphs must be inserted
# immediately before it (right after the "... oyuncu bir takımda bulunur."
# paragraph that currently precedes it).
$cokCok = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^\s*.OK-.OK\s*$") {
        $cokCok = $p
        break
    }
}

if ($cokCok -ne $null) {
    # Create a brand-new, empty paragraph right before "ÇOK-ÇOK" without
    # disturbing ÇOK-ÇOK's own paragraph mark/identity.
    $insertionRange = $cokCok.Range
    $insertionRange.Collapse(1)   # wdCollapseStart
    $insertionRange.InsertParagraphBefore()

    # Re-resolve "ÇOK-ÇOK" and grab the (now existing) empty paragraph that
    # sits right before it so we can fill it in with the real content.
    $cokCok = $null
    $prevPara = $null
    $all = @($d.Paragraphs)
    for ($i = 0; $i -lt $all.Count; $i++) {
        if ($all[$i].Range.Text -match "^\s*.OK-.OK\s*$") {
            $cokCok = $all[$i]
            $prevPara = $all[$i - 1]
            break
        }
    }

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">       4)</w:t></w:r><w:r><w:t xml:space="preserve">    Bir bilgisayarda birçok belge bulunabilir.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">               Birçok belge bir bilgisayarda bulunur.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $prevPara.Range.InsertXML($xml)
}
